# Regenerate merged AHB files
# - Rename header row labels: "_old" -> "_FV2304", "_new" -> "_FV2310"
# - Turn the data range into an Excel Table (ListObject) with an AutoFilter
# - Freeze the header row (top row) in the sheet view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename header row (row 1) cell values ---
# Columns A-J: "<Name>_old" -> "<Name>_FV2304"
$newHeaders2304 = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)
for ($i = 0; $i -lt $newHeaders2304.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $newHeaders2304[$i]
}

# Column K (11): "diff" is unchanged

# Columns L-U: "<Name>_new" -> "<Name>_FV2310"
$newHeaders2310 = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)
for ($i = 0; $i -lt $newHeaders2310.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $newHeaders2310[$i]
}

# --- 2) Convert the used range into an Excel Table with AutoFilter ---
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U60"), $null, 1, $null)
$tbl.Name = "Table1"

# --- 3) Freeze the top (header) row ---
$ws.Range("A2").Activate()
$excel.ActiveWindow.FreezePanes = $true
